# RESTORE_END v0.3 - Added passenger cars (fossil) - Plotting improvements
#
# The "passenger" flow values were stored in pkm (person-kilometres) but are
# being converted to Mpkm (million person-kilometres): every value in column
# G (rows 6-35) is divided by 1,000,000 and the unit label in column H is
# changed from "pkm" to "Mpkm".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 6; $r -le 35; $r++) {
    $oldValue = $ws.Range("G$r").Value2
    $ws.Range("G$r").Value = $oldValue / 1000000
    $ws.Range("H$r").Value = "Mpkm"
}

# Move the selection cursor (cosmetic, matches the saved view state).
[void]$ws.Range("J6").Select()
